$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2 DataFrames")

$ws.Range("B2").Value = 8.252643144094897
$ws.Range("C2").Value = 5.405017945584013
$ws.Range("F2").Value = 16.50528628818979
$ws.Range("G2").Value = 10.81003589116803

$ws.Range("B3").Value = 6.68588770720516
$ws.Range("C3").Value = 9.925574752198022
$ws.Range("F3").Value = 13.37177541441032
$ws.Range("G3").Value = 19.85114950439604

$ws.Range("B4").Value = 5.863923788494891
$ws.Range("C4").Value = 3.189183188888812
$ws.Range("F4").Value = 11.72784757698978
$ws.Range("G4").Value = 6.378366377777624

$ws.Range("B5").Value = 0.976589030820818
$ws.Range("C5").Value = 3.171506025267213
$ws.Range("F5").Value = 1.953178061641636
$ws.Range("G5").Value = 6.343012050534426

$ws.Range("B6").Value = 6.713361507575395
$ws.Range("C6").Value = 2.976258459102171
$ws.Range("F6").Value = 13.42672301515079
$ws.Range("G6").Value = 5.952516918204342

$ws.Range("B7").Value = 0.5099811919552821
$ws.Range("C7").Value = 0.01742846552266974
$ws.Range("F7").Value = 1.019962383910564
$ws.Range("G7").Value = 0.03485693104533949

$ws.Range("B8").Value = 2.891980130967887
$ws.Range("C8").Value = 8.653387117681449
$ws.Range("F8").Value = 5.783960261935775
$ws.Range("G8").Value = 17.3067742353629

$ws.Range("B9").Value = 9.495778467361443
$ws.Range("C9").Value = 9.931725249980042
$ws.Range("F9").Value = 18.99155693472289
$ws.Range("G9").Value = 19.86345049996008

$ws.Range("B10").Value = 9.735840553392082
$ws.Range("C10").Value = 6.241923081743693
$ws.Range("F10").Value = 19.47168110678416
$ws.Range("G10").Value = 12.48384616348739

$ws.Range("B11").Value = 4.223498224820437
$ws.Range("C11").Value = 8.485493146619246
$ws.Range("F11").Value = 8.446996449640874
$ws.Range("G11").Value = 16.97098629323849
